# Added 1 response from another council to test graphs
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Normalize columns K, M, S, U for existing rows 2-25 -----------------
# These columns currently hold numeric-looking answers ("1", "2", "3", "4")
# stored as text (shared strings). Re-writing them as real numbers keeps
# the displayed value identical but changes the underlying storage type,
# which also causes the now-unused "2" / "3" text entries to drop out of
# the shared-strings table automatically.
$numericCols = @(11, 13, 19, 21)  # K, M, S, U
for ($r = 2; $r -le 25; $r++) {
    foreach ($c in $numericCols) {
        $cell = $ws.Cells.Item($r, $c)
        $text = $cell.Value2
        $cell.Value = [double]$text
    }
}

# --- Add a brand new response row (row 26) --------------------------------
$ws.Cells.Item(26, 1).Value  = "Quarter 1"
$ws.Cells.Item(26, 2).Value  = 2
$ws.Cells.Item(26, 3).Value  = 0
$ws.Cells.Item(26, 4).Value  = 1
$ws.Cells.Item(26, 5).Value  = 1
$ws.Cells.Item(26, 6).Value  = 0
$ws.Cells.Item(26, 7).Value  = 1
$ws.Cells.Item(26, 8).Value  = 0
$ws.Cells.Item(26, 9).Value  = 0
$ws.Cells.Item(26, 10).Value = 0
$ws.Cells.Item(26, 11).Value = 3
$ws.Cells.Item(26, 12).Value = "-"
$ws.Cells.Item(26, 13).Value = 1
$ws.Cells.Item(26, 14).Value = "-"
$ws.Cells.Item(26, 15).Value = 4
$ws.Cells.Item(26, 16).Value = 1
$ws.Cells.Item(26, 17).Value = 2
$ws.Cells.Item(26, 18).Value = "-"
$ws.Cells.Item(26, 19).Value = 3
$ws.Cells.Item(26, 20).Value = "-"
$ws.Cells.Item(26, 21).Value = 4
$ws.Cells.Item(26, 22).Value = "-"
$ws.Cells.Item(26, 23).Value = "-"

# --- Update the view state (selection + scroll position) -----------------
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 16
$ws.Range("S28").Select()
